{"js": "// Replace each old three-digit-by-one-digit multiplication expression\n// with its corresponding new expression, matching the commit diff exactly.\nconst replacements = [\n  [\"811\u00d72=1622\", \"940\u00d78=7520\"],\n  [\"749\u00d74=2996\", \"399\u00d72=798\"],\n  [\"146\u00d79=1314\", \"958\u00d72=1916\"],\n  [\"832\u00d77=5824\", \"770\u00d73=2310\"],\n  [\"418\u00d75=2090\", \"512\u00d76=3072\"],\n  [\"253\u00d77=1771\", \"686\u00d74=2744\"],\n  [\"353\u00d75=1765\", \"454\u00d79=4086\"],\n  [\"788\u00d73=2364\", \"573\u00d79=5157\"],\n  [\"429\u00d77=3003\", \"649\u00d73=1947\"],\n  [\"616\u00d72=1232\", \"613\u00d79=5517\"],\n  [\"602\u00d79=5418\", \"814\u00d79=7326\"],\n  [\"168\u00d76=1008\", \"520\u00d74=2080\"],\n  [\"512\u00d79=4608\", \"926\u00d72=1852\"],\n  [\"664\u00d75=3320\", \"633\u00d73=1899\"],\n  [\"385\u00d74=1540\", \"713\u00d78=5704\"],\n  [\"673\u00d79=6057\", \"502\u00d79=4518\"],\n  [\"139\u00d79=1251\", \"120\u00d76=720\"],\n  [\"444\u00d75=2220\", \"452\u00d75=2260\"],\n  [\"882\u00d75=4410\", \"133\u00d78=1064\"],\n  [\"498\u00d77=3486\", \"239\u00d73=717\"],\n  [\"532\u00d78=4256\", \"648\u00d72=1296\"],\n  [\"682\u00d77=4774\", \"571\u00d77=3997\"],\n  [\"824\u00d77=5768\", \"477\u00d79=4293\"],\n  [\"389\u00d74=1556\", \"769\u00d75=3845\"],\n  [\"193\u00d78=1544\", \"145\u00d73=435\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each old three-digit-by-one-digit multiplication expression\n# with its corresponding new expression, matching the commit diff exactly.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"811\u00d72=1622\", \"940\u00d78=7520\"),\n  @(\"749\u00d74=2996\", \"399\u00d72=798\"),\n  @(\"146\u00d79=1314\", \"958\u00d72=1916\"),\n  @(\"832\u00d77=5824\", \"770\u00d73=2310\"),\n  @(\"418\u00d75=2090\", \"512\u00d76=3072\"),\n  @(\"253\u00d77=1771\", \"686\u00d74=2744\"),\n  @(\"353\u00d75=1765\", \"454\u00d79=4086\"),\n  @(\"788\u00d73=2364\", \"573\u00d79=5157\"),\n  @(\"429\u00d77=3003\", \"649\u00d73=1947\"),\n  @(\"616\u00d72=1232\", \"613\u00d79=5517\"),\n  @(\"602\u00d79=5418\", \"814\u00d79=7326\"),\n  @(\"168\u00d76=1008\", \"520\u00d74=2080\"),\n  @(\"512\u00d79=4608\", \"926\u00d72=1852\"),\n  @(\"664\u00d75=3320\", \"633\u00d73=1899\"),\n  @(\"385\u00d74=1540\", \"713\u00d78=5704\"),\n  @(\"673\u00d79=6057\", \"502\u00d79=4518\"),\n  @(\"139\u00d79=1251\", \"120\u00d76=720\"),\n  @(\"444\u00d75=2220\", \"452\u00d75=2260\"),\n  @(\"882\u00d75=4410\", \"133\u00d78=1064\"),\n  @(\"498\u00d77=3486\", \"239\u00d73=717\"),\n  @(\"532\u00d78=4256\", \"648\u00d72=1296\"),\n  @(\"682\u00d77=4774\", \"571\u00d77=3997\"),\n  @(\"824\u00d77=5768\", \"477\u00d79=4293\"),\n  @(\"389\u00d74=1556\", \"769\u00d75=3845\"),\n  @(\"193\u00d78=1544\", \"145\u00d73=435\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
